$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '40.018.01'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.63%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.238.77'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -3.74%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '294.47'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -4.46%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '86.14'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +3.55%  '

$ws.Range('E7').Value = '  -1.73%  '

$ws.Range('E8').Value = '  +0.02%  '

$ws.Range('E9').Value = '  -1.11%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0797'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.06%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '30.49'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +3.35%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '47.45'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -9.39%  '

$ws.Range('E13').Value = '  -2.09%  '

$ws.Range('E14').Value = '  +0.82%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.583.90'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -3.77%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.21'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -2.51%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.234.39'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -4.91%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.725'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -3.04%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '39.927.98'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.55%  '

$ws.Range('E20').Value = '  +0.23%  '

$ws.Range('E21').Value = '  -3.26%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.71'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +3.00%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '65.54'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -3.18%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '234.63'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.63%  '

$ws.Range('E25').Value = '  +0.00%  '

$ws.Range('E26').Value = '  -3.10%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.85'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.83%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '23.03'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.84%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.21'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.54%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '9.24'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.89%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '33.55'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.04%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '155.13'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.85%  '

$ws.Range('E33').Value = '  -0.23%  '

$ws.Range('E34').Value = '  -3.29%  '

$ws.Range('E35').Value = '  +0.34%  '

$ws.Range('E36').Value = '  -4.24%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '16.49'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +8.01%  '

$ws.Range('E38').Value = '  -0.67%  '

$ws.Range('E39').Value = '  +2.23%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.69'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -1.43%  '

$ws.Range('E41').Value = '  -0.67%  '

$ws.Range('E42').Value = '  +2.21%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.954.24'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.63%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.18'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -3.10%  '

$ws.Range('E45').Value = '  +3.99%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '9.52'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.51%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '16.33'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -4.72%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.61'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.88%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.455.48'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -3.95%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '70.96'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.24%  '

$ws.Range('E51').Value = '  +9.22%  '
